$d = $word.ActiveDocument

# The table logging "Fecha" / "Objetivo realizado" entries needs a new
# row appended at the bottom: 11/08/2020 -> "Prototipo ciudad html"
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "11/08/2020"
$newRow.Cells.Item(2).Range.Text = "Prototipo ciudad html"
